$wb = $excel.ActiveWorkbook

# Rename sheets (task order titles with updated timestamps)
$wb.Worksheets.Item(1).Name = "GNG_TO-16504778558065643"
$wb.Worksheets.Item(2).Name = "NB_TO-16504778583995283"
$wb.Worksheets.Item(3).Name = "RS_TO-16504778584015288"
$wb.Worksheets.Item(4).Name = "TOL_TO-16504778584485295"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504778585115669"

# Sheet1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504778557784436.csv"
$ws1.Range("B3").Value = "GNG_stims-1650477855789528.csv"
$ws1.Range("B4").Value = "go_stims-16504778557905524.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778558055615.csv"

# Sheet2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16504778572475295.csv"
$ws2.Range("B3").Value = "ZB-match_2-1650477856218532.csv"
$ws2.Range("B4").Value = "OB-16504778567615316.csv"
$ws2.Range("B5").Value = "ZB-match_4-1650477855896531.csv"
$ws2.Range("B6").Value = "TB-16504778582325628.csv"
$ws2.Range("B7").Value = "ZB-match_0-16504778559355304.csv"
$ws2.Range("B8").Value = "TB-1650477858380526.csv"
$ws2.Range("B9").Value = "OB-16504778566105583.csv"
$ws2.Range("B10").Value = "OB-16504778570135624.csv"

# Sheet4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504778584155633.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778584035282.csv"
$ws4.Range("B4").Value = "MM_stims-16504778584315627.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778584155633.csv"
$ws4.Range("B6").Value = "MM_stims-16504778584475634.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778584325292.csv"

# Sheet5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16504778584635272.csv"
$ws5.Range("B3").Value = "vSAT_stims-16504778584955637.csv"
$ws5.Range("B4").Value = "vSAT_stims-1650477858479564.csv"
$ws5.Range("B5").Value = "SAT_stims-16504778584515293.csv"
